# Actualización automática desde tarea programada
# Adds a new measurement row (row 14) to the sheet and nudges the
# timestamp of the previous last row (row 13) as produced by the
# scheduled data-collection task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- tiny correction to the existing last row's timestamp -----------------
$ws.Range("A13").Value = 45868.50022140046

# --- append the new row of sensor readings ---------------------------------
$newRow = 14

$ws.Range("A" + $newRow).NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("A" + $newRow).Value = 45868.54185007951
$ws.Range("B" + $newRow).Value = 2025
$ws.Range("C" + $newRow).Value = 31
$ws.Range("D" + $newRow).Value = 23.44
$ws.Range("E" + $newRow).Value = 66.22
$ws.Range("F" + $newRow).Value = 613.48
$ws.Range("G" + $newRow).Value = 10.78
$ws.Range("H" + $newRow).Value = "ESE"
$ws.Range("I" + $newRow).Value = 0
$ws.Range("J" + $newRow).Value = "13:00:15"
